# Analysis.xlsx edit:
#  - Rename Sheet1 -> Data, Sheet2 -> Calculations
#  - Make "Data" the active/selected sheet (was "Sheet2"/Calculations)
#  - Update the "Data" sheet view: scroll to top-left A5, zoom 75%,
#    selection moves from Q1 to O1
#  - Update the "Calculations" sheet view (no longer active tab):
#    selection moves from C8 to C7

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item(1)
$wsCalc = $wb.Worksheets.Item(2)

# --- Rename the worksheets ---
$wsData.Name = "Data"
$wsCalc.Name = "Calculations"

# --- Update "Calculations" sheet view/selection first (it will stop being
#     the active tab once "Data" is (re)activated below) ---
$wsCalc.Activate()
$wsCalc.Range("C7").Select() | Out-Null

# --- Update "Data" sheet view/selection and make it the active tab ---
$wsData.Activate()
$wData = $excel.ActiveWindow
$wsData.Range("O1").Select() | Out-Null
$wData.Zoom = 75
$wData.ScrollRow = 5
$wData.ScrollColumn = 1

# Ensure "Data" ends up as the active sheet (tabSelected) of the workbook.
$wsData.Activate()
